$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Moorings")
$ws2 = $wb.Worksheets.Item("Asset_Cal_Info")

# --- Asset_Cal_Info sheet: correct instrument reference designators ---
# (glider CP05MOAS-GL001 was actually CP05MOAS-GL336) and deployment number
$ws2.Range("A2:A5").Value = "CP05MOAS-GL336-01-ADCPAM000"
$ws2.Range("A7:A10").Value = "CP05MOAS-GL336-02-FLORTM000"
$ws2.Range("A12").Value = "CP05MOAS-GL336-03-CTDGVM000"
$ws2.Range("A14").Value = "CP05MOAS-GL336-04-DOSTAM000"
$ws2.Range("A16").Value = "CP05MOAS-GL336-05-PARADM000"
$ws2.Range("A18").Value = "CP05MOAS-GL336-00-ENG000000"

$ws2.Range("C2").Value = 1
$ws2.Range("C3").Value = 1
$ws2.Range("C4").Value = 1
$ws2.Range("C5").Value = 1
$ws2.Range("C7").Value = 1
$ws2.Range("C8").Value = 1
$ws2.Range("C9").Value = 1
$ws2.Range("C10").Value = 1
$ws2.Range("C12").Value = 1
$ws2.Range("C14").Value = 1
$ws2.Range("C16").Value = 1
$ws2.Range("C18").Value = 1

# --- Moorings sheet: correct Ref Des and deployment number ---
$ws1.Range("A2").Value = "CP05MOAS-GL336"
$ws1.Range("C2").Value = 1

# --- restore the active-tab / selection state that was captured at save time ---
$ws2.Activate()
$ws2.Range("C20").Select()

$ws1.Activate()
$ws1.Range("B15").Select()
